$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '46.807.52'
$ws.Range('E2').Value = '  +4.89%  '
$ws.Range('D3').Value = '2.341.21'
$ws.Range('E3').Value = '  +4.28%  '
$ws.Range('E4').Value = '  -0.62%  '
$ws.Range('D5').Value = "'305.34"
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = "'97.13"
$ws.Range('E6').Value = '  +2.72%  '
$ws.Range('D7').Value = "'0.576"
$ws.Range('E7').Value = '  +1.22%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').Value = "'0.534"
$ws.Range('E9').Value = '  +3.67%  '
$ws.Range('D10').Value = "'35.65"
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').Value = "'7.39"
$ws.Range('E12').Value = '  +2.85%  '
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').Value = '2.698.45'
$ws.Range('E14').Value = '  +4.23%  '
$ws.Range('D15').Value = '2.342.05'
$ws.Range('E15').Value = '  -1.51%  '
$ws.Range('D16').Value = "'14.14"
$ws.Range('E16').Value = '  +4.46%  '
$ws.Range('D17').Value = "'0.827"
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').Value = '46.682.96'
$ws.Range('E18').Value = '  +5.16%  '
$ws.Range('D19').Value = "'13.60"
$ws.Range('E19').Value = '  +15.97%  '
$ws.Range('D20').Value = '0.0₃0946'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').Value = "'6.19"
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('E22').Value = '  +2.42%  '
$ws.Range('D23').Value = "'245.13"
$ws.Range('E23').Value = '  +3.21%  '
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('D27').Value = "'41.45"
$ws.Range('E27').Value = '  +12.27%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').Value = "'9.85"
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('D31').Value = "'5.74"
$ws.Range('E31').Value = '  -1.75%  '
$ws.Range('D32').Value = "'152.48"
$ws.Range('E32').Value = '  +3.23%  '
$ws.Range('D33').Value = "'0.0813"
$ws.Range('E33').Value = '  +4.20%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').Value = "'3.17"
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('E36').Value = '  +1.65%  '
$ws.Range('D37').Value = "'0.118"
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('D38').Value = "'1.81"
$ws.Range('E38').Value = '  -2.19%  '
$ws.Range('D39').Value = "'4.02"
$ws.Range('E39').Value = '  +7.10%  '
$ws.Range('D40').Value = "'0.0314"
$ws.Range('E40').Value = '  +5.31%  '
$ws.Range('D41').Value = "'3.38"
$ws.Range('E41').Value = '  +1.34%  '
$ws.Range('D42').Value = "'13.71"
$ws.Range('E42').Value = '  -9.50%  '
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').Value = "'1.97"
$ws.Range('E44').Value = '  +11.98%  '
$ws.Range('D45').Value = '1.839.52'
$ws.Range('E45').Value = '  +1.68%  '
$ws.Range('D46').Value = "'0.196"
$ws.Range('E46').Value = '  +5.21%  '
$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').Value = "'73.77"
$ws.Range('E47').Value = '  +7.23%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').Value = "'80.73"
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('D49').Value = "'4.93"
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').Value = "'55.01"
$ws.Range('E51').Value = '  +2.07%  '
